$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3.8
$ws.Range("I8").Value = 3.8
$ws.Range("K8").Value = 11.4
$ws.Range("M8").Value = 127.6

$ws.Range("H70").Value = 4005.5715
$ws.Range("I70").Value = 3499.5
$ws.Range("K70").Value = 10498.5
$ws.Range("M70").Value = -10228.5

$ws.Range("H73").Value = 4005.5715
$ws.Range("I73").Value = 3499.5
$ws.Range("K73").Value = 10498.5
$ws.Range("M73").Value = -9562.5

$ws.Range("H74").Value = 102559.8
$ws.Range("I74").Value = 3099.6667
$ws.Range("J74").Value = 251750
$ws.Range("K74").Value = 3099.6667
$ws.Range("L74").Value = 251750
$ws.Range("M74").Value = -2163.6667
$ws.Range("N74").Value = -253622

$ws.Range("H77").Value = 102559.8
$ws.Range("I77").Value = 3099.6667
$ws.Range("J77").Value = 251750
$ws.Range("K77").Value = 15498.3335
$ws.Range("L77").Value = 1258750
$ws.Range("M77").Value = -10818.3335
$ws.Range("N77").Value = -1268110

$ws.Range("H80").Value = 971.94446
$ws.Range("I80").Value = 940
$ws.Range("J80").Value = 984.2308
$ws.Range("K80").Value = 2820
$ws.Range("L80").Value = 2952.6924
$ws.Range("M80").Value = -1822
$ws.Range("N80").Value = -4948.6924

$ws.Range("H83").Value = 971.94446
$ws.Range("I83").Value = 940
$ws.Range("J83").Value = 984.2308
$ws.Range("K83").Value = 8460
$ws.Range("L83").Value = 8858.0772
$ws.Range("M83").Value = -3468
$ws.Range("N83").Value = -18842.0772

$ws.Range("H99").Value = 1115
$ws.Range("I99").Value = 633
$ws.Range("J99").Value = 1597
$ws.Range("K99").Value = 1899
$ws.Range("L99").Value = 4791
$ws.Range("M99").Value = -401
$ws.Range("N99").Value = -7787

$ws.Range("H101").Value = 516
$ws.Range("J101").Value = 582.5
$ws.Range("L101").Value = 1747.5
$ws.Range("N101").Value = -4991.5

$ws.Range("H137").Value = 2133.3572
$ws.Range("I137").Value = 1670.6316
$ws.Range("K137").Value = 5011.8948
$ws.Range("M137").Value = -2461.8948

$ws.Range("H138").Value = 2213.4119
$ws.Range("I138").Value = 1096.8572
$ws.Range("K138").Value = 3290.5716
$ws.Range("M138").Value = 1849.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2323.16
$ws.Range("I32").Value = 2003.2916
$ws.Range("K32").Value = 2003.2916
$ws.Range("M32").Value = -1716.2916

$ws.Range("H45").Value = 1369.5
$ws.Range("I45").Value = 949.1667
$ws.Range("K45").Value = 949.1667
$ws.Range("M45").Value = -572.1667

$ws.Range("H61").Value = 1575.2
$ws.Range("I61").Value = 1575.2
$ws.Range("K61").Value = 1575.2
$ws.Range("M61").Value = -1363.2

$ws.Range("H122").Value = 3082.5
$ws.Range("I122").Value = 1499.6666
$ws.Range("K122").Value = 4498.9998
$ws.Range("M122").Value = -2048.9998

$ws.Range("H136").Value = 1575.2
$ws.Range("I136").Value = 1575.2
$ws.Range("K136").Value = 4725.6
$ws.Range("M136").Value = -2175.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5499.5
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 5499.5
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3507.5
$ws.Range("I58").Value = 1609
$ws.Range("J58").Value = 5948.4287
$ws.Range("K58").Value = 1609
$ws.Range("L58").Value = 5948.4287
$ws.Range("M58").Value = -1406
$ws.Range("N58").Value = -6354.4287

$ws.Range("H99").Value = 4213.737
$ws.Range("I99").Value = 3625.3333
$ws.Range("J99").Value = 4743.3
$ws.Range("K99").Value = 3625.3333
$ws.Range("L99").Value = 4743.3
$ws.Range("M99").Value = -2127.3333
$ws.Range("N99").Value = -7739.3

$ws.Range("H107").Value = 433.16666
$ws.Range("I107").Value = 433.16666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 433.16666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1486.83334
$ws.Range("N107").ClearContents()

$ws.Range("H126").Value = 4213.737
$ws.Range("I126").Value = 3625.3333
$ws.Range("J126").Value = 4743.3
$ws.Range("K126").Value = 10875.9999
$ws.Range("L126").Value = 14229.9
$ws.Range("M126").Value = -8405.999899999999
$ws.Range("N126").Value = -19169.9

$ws.Range("H136").Value = 3507.5
$ws.Range("I136").Value = 1609
$ws.Range("J136").Value = 5948.4287
$ws.Range("K136").Value = 4827
$ws.Range("L136").Value = 17845.2861
$ws.Range("M136").Value = -2277
$ws.Range("N136").Value = -22945.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 522.2222
$ws.Range("I107").Value = 419.6
$ws.Range("K107").Value = 1258.8
$ws.Range("M107").Value = 661.1999999999998

$ws.Range("H121").Value = 672.0909
$ws.Range("I121").Value = 470.57144
$ws.Range("K121").Value = 1411.71432
$ws.Range("M121").Value = -101.71432

$ws.Range("H141").Value = 1959.7778
$ws.Range("I141").Value = 1959.7778
$ws.Range("K141").Value = 5879.3334
$ws.Range("M141").Value = -699.3334000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 40000
$ws.Range("J47").Value = 40000
$ws.Range("L47").Value = 40000
$ws.Range("N47").Value = -41136

$ws.Range("H113").Value = 3499.8
$ws.Range("J113").Value = 3899.75
$ws.Range("L113").Value = 3899.75
$ws.Range("N113").Value = -8239.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 146.5
$ws.Range("I55").Value = 89.5
$ws.Range("K55").Value = 89.5
$ws.Range("M55").Value = 83.5

$ws.Range("H93").Value = 3028.2856
$ws.Range("I93").Value = 3033
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 3033
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1785
$ws.Range("N93").Value = -5496

$ws.Range("H122").Value = 3590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 25425.75
$ws.Range("J4").Value = 25425.75
$ws.Range("L4").Value = 25425.75
$ws.Range("N4").Value = -25651.75

$ws.Range("H8").Value = 1001.5
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280

$ws.Range("H11").Value = 41731668
$ws.Range("J11").Value = 97500
$ws.Range("L11").Value = 97500
$ws.Range("N11").Value = -97784

$ws.Range("H62").Value = 3447.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3447.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3447.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4695.5

$ws.Range("H65").Value = 3447.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3447.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 17237.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -23477.5

$ws.Range("H81").Value = 1204.5
$ws.Range("I81").Value = 1045.4
$ws.Range("K81").Value = 2090.8
$ws.Range("M81").Value = -1029.8

$ws.Range("H84").Value = 1204.5
$ws.Range("I84").Value = 1045.4
$ws.Range("K84").Value = 10454
$ws.Range("M84").Value = -5150

$ws.Range("H113").Value = 290.6
$ws.Range("I113").Value = 248.5
$ws.Range("K113").Value = 745.5
$ws.Range("M113").Value = 1424.5
